$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 11:07"

# Row 7: India -> India
$ws.Cells.Item(7, 2).Value = 355060
$ws.Cells.Item(7, 3).Value = 899
$ws.Cells.Item(7, 4).Value = 187718
$ws.Cells.Item(7, 5).Value = 155420
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 8).Value = 11922

# Row 21: Banglades -> Banglades
$ws.Cells.Item(21, 2).Value = 98489
$ws.Cells.Item(21, 3).Value = 4008
$ws.Cells.Item(21, 4).Value = 38189
$ws.Cells.Item(21, 5).Value = 58995
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 43
$ws.Cells.Item(21, 8).Value = 1305

# Row 40: Polonia -> Polonia
$ws.Cells.Item(40, 2).Value = 30701
$ws.Cells.Item(40, 3).Value = 506
$ws.Cells.Item(40, 4).Value = 14921
$ws.Cells.Item(40, 5).Value = 14494
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 14
$ws.Cells.Item(40, 8).Value = 1286

# Row 41: Afganistan -> Filipinas
$ws.Cells.Item(41, 1).Value = "Filipinas"
$ws.Cells.Item(41, 2).Value = 27238
$ws.Cells.Item(41, 3).Value = 457
$ws.Cells.Item(41, 4).Value = 6820
$ws.Cells.Item(41, 5).Value = 19310
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 5
$ws.Cells.Item(41, 8).Value = 1108

# Row 42: Filipinas -> Afganistan
$ws.Cells.Item(42, 1).Value = "Afganistan"
$ws.Cells.Item(42, 2).Value = 26874
$ws.Cells.Item(42, 3).Value = 564
$ws.Cells.Item(42, 4).Value = 6158
$ws.Cells.Item(42, 5).Value = 20212
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 13
$ws.Cells.Item(42, 8).Value = 504

# Row 43: Irlanda -> Oman
$ws.Cells.Item(43, 1).Value = "Oman"
$ws.Cells.Item(43, 2).Value = 26079
$ws.Cells.Item(43, 3).Value = 810
$ws.Cells.Item(43, 4).Value = 11797
$ws.Cells.Item(43, 5).Value = 14166
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = 116

# Row 44: Oman -> Irlanda
$ws.Cells.Item(44, 1).Value = "Irlanda"
$ws.Cells.Item(44, 2).Value = 25334
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 22698
$ws.Cells.Item(44, 5).Value = 927
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 1709

# Row 56: Kazajistan -> Kazajistan
$ws.Cells.Item(56, 2).Value = 15542
$ws.Cells.Item(56, 3).Value = 350
$ws.Cells.Item(56, 4).Value = 9716
$ws.Cells.Item(56, 5).Value = 5738
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 88

# Row 58: Moldavia -> Moldavia
$ws.Cells.Item(58, 2).Value = 12254
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 7077
$ws.Cells.Item(58, 5).Value = 4750
$ws.Cells.Item(58, 6).Value = 0
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 427

# Row 87: El Salvador -> El Salvador
$ws.Cells.Item(87, 2).Value = 4066
$ws.Cells.Item(87, 3).Value = 125
$ws.Cells.Item(87, 4).Value = 2137
$ws.Cells.Item(87, 5).Value = 1851
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 78

# Row 109: Sudan del Sur -> Lituania
$ws.Cells.Item(109, 1).Value = "Lituania"
$ws.Cells.Item(109, 2).Value = 1778
$ws.Cells.Item(109, 3).Value = 2
$ws.Cells.Item(109, 4).Value = 1447
$ws.Cells.Item(109, 5).Value = 255
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 76

# Row 110: Lituania -> Sudan del Sur
$ws.Cells.Item(110, 1).Value = "Sudan del Sur"
$ws.Cells.Item(110, 2).Value = 1776
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(110, 4).Value = 58
$ws.Cells.Item(110, 5).Value = 1688
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 0
$ws.Cells.Item(110, 8).Value = 30

# Row 113: Eslovaquia -> Eslovaquia
$ws.Cells.Item(113, 2).Value = 1561
$ws.Cells.Item(113, 3).Value = 9
$ws.Cells.Item(113, 4).Value = 1437
$ws.Cells.Item(113, 5).Value = 96
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 28

# Row 137: Uganda -> Uganda
$ws.Cells.Item(137, 2).Value = 732
$ws.Cells.Item(137, 3).Value = 8
$ws.Cells.Item(137, 4).Value = 409
$ws.Cells.Item(137, 5).Value = 323
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 0

# Row 145: Malaui -> Benin
$ws.Cells.Item(145, 1).Value = "Benin"
$ws.Cells.Item(145, 2).Value = 572
$ws.Cells.Item(145, 3).Value = 40
$ws.Cells.Item(145, 4).Value = 237
$ws.Cells.Item(145, 5).Value = 326
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 9

# Row 146: Togo -> Malaui
$ws.Cells.Item(146, 1).Value = "Malaui"
$ws.Cells.Item(146, 2).Value = 564
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 73
$ws.Cells.Item(146, 5).Value = 485
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 6

# Row 147: Benin -> Togo
$ws.Cells.Item(147, 1).Value = "Togo"
$ws.Cells.Item(147, 2).Value = 537
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 344
$ws.Cells.Item(147, 5).Value = 180
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 13

# Row 206: Groenlandia -> Islas Malvinas
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 207: Islas Malvinas -> Groenlandia
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 210: Seychelles -> Montserrat
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 1

# Row 211: Montserrat -> Seychelles
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0

# Row 213: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

